$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1): "_old" columns -> "_FV2210", "_new" -> "_FV2304"
#    ("diff" in K1 is left untouched.)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# ---------------------------------------------------------------------------
# 2) Turn A1:U73 into an Excel Table ("Table1") without disturbing the
#    existing header-row formatting/styles.xml. Excel's "insert table"
#    normally freezes a copy of any pre-existing header formatting into a
#    new dxf (headerRowDxfId) - stash + restore the formatting via
#    copy/paste so no new style entries get created.
# ---------------------------------------------------------------------------
$ws.Range("A1:U1").Copy()
$ws.Range("A1000:U1000").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1:U1").ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U73"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$ws.Range("A1000:U1000").Copy()
$ws.Range("A1:U1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1000:U1000").Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
